$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 12-35 (their fields get folded into rows 2-11)
$ws.Range("A12:A35").EntireRow.Delete() | Out-Null

# Rewrite rows 2-11 with the consolidated "(name, [fields...])" text
$ws.Range("A2").Value = '(''Bloodstained Mire'', [''Land'', ''{T}, Pay 1 life, Sacrifice Bloodstained Mire: Search your library for a Swamp or Mountain card, put it onto the battlefield, then shuffle your library.''])'
$ws.Range("A3").Value = '(''Burning Wish'', [''{1}{R}'', ''Sorcery'', ''You may reveal a sorcery card you own from outside the game and put it into your hand. Exile Burning Wish.''])'
$ws.Range("A4").Value = '(''Dark Ritual'', [''{B}'', ''Instant'', ''Add {B}{B}{B}.''])'
$ws.Range("A5").Value = '(''Flooded Strand'', [''Land'', ''{T}, Pay 1 life, Sacrifice Flooded Strand: Search your library for a Plains or Island card, put it onto the battlefield, then shuffle your library.''])'
$ws.Range("A6").Value = '(''Maze of Ith'', [''Land'', ''{T}: Untap target attacking creature. Prevent all combat damage that would be dealt to and dealt by that creature this turn.''])'
$ws.Range("A7").Value = '(''Polluted Delta'', [''Land'', ''{T}, Pay 1 life, Sacrifice Polluted Delta: Search your library for an Island or Swamp card, put it onto the battlefield, then shuffle your library.''])'
$ws.Range("A8").Value = '(''Stifle'', [''{U}'', ''Instant'', ''Counter target activated or triggered ability. (Mana abilities can’t be targeted.)''])'
$ws.Range("A9").Value = '(''Survival of the Fittest'', [''{1}{G}'', ''Enchantment'', ''{G}, Discard a creature card: Search your library for a creature card, reveal that card, and put it into your hand. Then shuffle your library.''])'
$ws.Range("A10").Value = '(''Windswept Heath'', [''Land'', ''{T}, Pay 1 life, Sacrifice Windswept Heath: Search your library for a Forest or Plains card, put it onto the battlefield, then shuffle your library.''])'
$ws.Range("A11").Value = '(''Wooded Foothills'', [''Land'', ''{T}, Pay 1 life, Sacrifice Wooded Foothills: Search your library for a Mountain or Forest card, put it onto the battlefield, then shuffle your library.''])'
